$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 598.8889
$ws.Range("I12").Value = 488
$ws.Range("J12").Value = 737.5
$ws.Range("K12").Value = 488
$ws.Range("L12").Value = 737.5
$ws.Range("M12").Value = -318
$ws.Range("N12").Value = -1077.5

$ws.Range("H92").Value = 145.05882
$ws.Range("I92").Value = 137.73334
$ws.Range("K92").Value = 137.73334
$ws.Range("M92").Value = 1110.26666

$ws.Range("H103").Value = 5552.6924
$ws.Range("I103").Value = 845
$ws.Range("K103").Value = 2535
$ws.Range("M103").Value = -1949

$ws.Range("H138").Value = 2101.3572
$ws.Range("I138").Value = 1201.5834
$ws.Range("J138").Value = 7500
$ws.Range("K138").Value = 3604.7502
$ws.Range("L138").Value = 22500
$ws.Range("M138").Value = 1535.2498
$ws.Range("N138").Value = -32780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3570.524
$ws.Range("I74").Value = 3271.7222
$ws.Range("J74").Value = 5363.3335
$ws.Range("K74").Value = 3271.7222
$ws.Range("L74").Value = 5363.3335
$ws.Range("M74").Value = -2397.7222
$ws.Range("N74").Value = -7111.3335

$ws.Range("H77").Value = 3570.524
$ws.Range("I77").Value = 3271.7222
$ws.Range("J77").Value = 5363.3335
$ws.Range("K77").Value = 16358.611
$ws.Range("L77").Value = 26816.6675
$ws.Range("M77").Value = -11990.611
$ws.Range("N77").Value = -35552.6675

$ws.Range("H88").Value = 895.3
$ws.Range("I88").Value = 419.4
$ws.Range("J88").Value = 1371.2
$ws.Range("K88").Value = 419.4
$ws.Range("L88").Value = 1371.2
$ws.Range("M88").Value = -13.39999999999998
$ws.Range("N88").Value = -2183.2

$ws.Range("H91").Value = 895.3
$ws.Range("I91").Value = 419.4
$ws.Range("J91").Value = 1371.2
$ws.Range("K91").Value = 419.4
$ws.Range("L91").Value = 1371.2
$ws.Range("M91").Value = 984.6
$ws.Range("N91").Value = -4179.2

$ws.Range("H97").Value = 2404.4285
$ws.Range("I97").Value = 2535
$ws.Range("J97").Value = 2230.3333
$ws.Range("K97").Value = 2535
$ws.Range("L97").Value = 2230.3333
$ws.Range("M97").Value = -2039
$ws.Range("N97").Value = -3222.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 629.6667
$ws.Range("J94").Value = 589
$ws.Range("L94").Value = 589
$ws.Range("N94").Value = -1491

$ws.Range("H99").Value = 83334240
$ws.Range("I99").Value = 100000870
$ws.Range("K99").Value = 100000870
$ws.Range("M99").Value = -99999372

$ws.Range("H107").Value = 25005080
$ws.Range("I107").Value = 55557510
$ws.Range("K107").Value = 55557510
$ws.Range("M107").Value = -55555590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3773.9688
$ws.Range("I31").Value = 2397.25
$ws.Range("K31").Value = 2397.25
$ws.Range("M31").Value = -2102.25

$ws.Range("H34").Value = 3773.9688
$ws.Range("I34").Value = 2397.25
$ws.Range("K34").Value = 2397.25
$ws.Range("M34").Value = -2195.25

$ws.Range("H107").Value = 857.2222
$ws.Range("I107").Value = 240
$ws.Range("K107").Value = 240
$ws.Range("M107").Value = 1680

$ws.Range("H134").Value = 3995.4443
$ws.Range("I134").Value = 2321.3333
$ws.Range("K134").Value = 6963.999899999999
$ws.Range("M134").Value = -4428.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1488.6364
$ws.Range("I57").Value = 1337.5
$ws.Range("K57").Value = 4012.5
$ws.Range("M57").Value = -3453.5

$ws.Range("H104").Value = 8944.333000000001
$ws.Range("J104").Value = 9937.375
$ws.Range("L104").Value = 29812.125
$ws.Range("N104").Value = -35054.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 3978.3333
$ws.Range("I44").Value = 4000
$ws.Range("J44").Value = 3935
$ws.Range("K44").Value = 4000
$ws.Range("L44").Value = 3935
$ws.Range("M44").Value = -3404
$ws.Range("N44").Value = -5127

$ws.Range("H80").Value = 4575
$ws.Range("I80").Value = 3600
$ws.Range("K80").Value = 3600
$ws.Range("M80").Value = -2602

$ws.Range("H83").Value = 4575
$ws.Range("I83").Value = 3600
$ws.Range("K83").Value = 18000
$ws.Range("M83").Value = -13008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2543.625
$ws.Range("I46").Value = 698.8
$ws.Range("J46").Value = 3382.182
$ws.Range("K46").Value = 698.8
$ws.Range("L46").Value = 3382.182
$ws.Range("M46").Value = -510.8
$ws.Range("N46").Value = -3758.182

$ws.Range("H61").Value = 62503692
$ws.Range("I61").Value = 125001144
$ws.Range("J61").Value = 6236.75
$ws.Range("K61").Value = 125001144
$ws.Range("L61").Value = 6236.75
$ws.Range("M61").Value = -125000942
$ws.Range("N61").Value = -6640.75

$ws.Range("H82").Value = 3399.9167
$ws.Range("I82").Value = 1008.3333
$ws.Range("K82").Value = 1008.3333
$ws.Range("M82").Value = -647.3333

$ws.Range("H85").Value = 3399.9167
$ws.Range("I85").Value = 1008.3333
$ws.Range("K85").Value = 1008.3333
$ws.Range("M85").Value = 239.6667

$ws.Range("H93").Value = 2959
$ws.Range("I93").Value = 2959
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2959
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1711
$ws.Range("N93").ClearContents()

$ws.Range("H113").Value = 62503692
$ws.Range("I113").Value = 125001144
$ws.Range("J113").Value = 6236.75
$ws.Range("K113").Value = 125001144
$ws.Range("L113").Value = 6236.75
$ws.Range("M113").Value = -124998974
$ws.Range("N113").Value = -10576.75

$ws.Range("H120").Value = 57000
$ws.Range("J120").Value = 57000
$ws.Range("L120").Value = 57000
$ws.Range("N120").Value = -66676

$ws.Range("H132").Value = 2976.2666
$ws.Range("I132").Value = 1604.5
$ws.Range("K132").Value = 4813.5
$ws.Range("M132").Value = -2283.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 4043000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 4043000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 4043000
$ws.Range("N54").Value = -4044040
$ws.Range("M54").ClearContents()

$ws.Range("H100").Value = 1371.2727
$ws.Range("I100").Value = 1371.2727
$ws.Range("K100").Value = 2742.5454
$ws.Range("M100").Value = -2201.5454

$ws.Range("H132").Value = 2153
$ws.Range("J132").Value = 3499.5
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -15558.5
